# Integration von profilgesteuerten Lasten
# Updates the line-length (D) and cos-phi/F values for the load profile rows,
# re-derives the shared "=15*D" formula results (handled automatically by the
# recalculation engine), applies a thousands-separator number format to the
# rows whose new length values are no longer "round" numbers, and moves the
# active selection to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 58
$ws.Range("F2").Value = 6.33245608

# Row 3
$ws.Range("D3").Value = 72
$ws.Range("F3").Value = 5.9615434240000003

# Row 4
$ws.Range("D4").Value = 304
$ws.Range("D4").NumberFormat = "#,##0"
$ws.Range("F4").Value = 5.1429563480000002

# Row 5
$ws.Range("D5").Value = 2.7
$ws.Range("F5").Value = 5.8600341299999998

# Row 6
$ws.Range("D6").Value = 101
$ws.Range("D6").NumberFormat = "#,##0"
$ws.Range("F6").Value = 5.371219601

# Row 7
$ws.Range("D7").Value = 101
$ws.Range("D7").NumberFormat = "#,##0"
$ws.Range("F7").Value = 4.860041152

# Row 8
$ws.Range("D8").Value = 540
$ws.Range("D8").NumberFormat = "#,##0"
$ws.Range("F8").Value = 4.6324939289999998

# Row 9
$ws.Range("D9").Value = 200
$ws.Range("D9").NumberFormat = "#,##0"
$ws.Range("F9").Value = 5.4378304499999999

# Row 10
$ws.Range("D10").Value = 200
$ws.Range("D10").NumberFormat = "#,##0"
$ws.Range("F10").Value = 6.2177166230000003

# Move the active selection to D8 (matches the author's last edit position)
[void]$ws.Range("D8").Select()
